# Rename the QR-code images (and matching URLs) on the "Access this page"
# and "Access this presentation" slides so that they point at the Overview
# deck/page instead of the Principles one.

function Update-UrlParagraph {
    param(
        $Shape,
        [int]$ParagraphIndex,
        [string]$NewText
    )

    $tr = $Shape.TextFrame.TextRange
    $para = $tr.Paragraphs($ParagraphIndex)
    # Replace the whole paragraph's character span in a single operation so
    # the existing run (and its formatting / hyperlink) is reused instead of
    # being split into a "common prefix" run + a "changed suffix" run.
    $chars = $para.Characters(1, $para.Length)
    $chars.Text = $NewText
}

$p = $ppt.ActivePresentation

# --- Slide 29: "Access this page" (HTML link + QR code) -------------------
$s29 = $p.Slides.Item(29)

$textShape29 = $s29.Shapes.Item("Text Placeholder 3")
Update-UrlParagraph $textShape29 2 "https://sbn-psi.github.io/dmsp/LDDTesting/LDDTestingOverview"

$picShape29 = $s29.Shapes.Item("Picture 1")
$picShape29.AlternativeText = "fig:  images/qr/overview_page.png"

# --- Slide 30: "Access this presentation" (PPTX link + QR code) -----------
$s30 = $p.Slides.Item(30)

$textShape30 = $s30.Shapes.Item("Text Placeholder 3")
Update-UrlParagraph $textShape30 2 "https://github.com/sbn-psi/dmsp/raw/main/LDDTesting/stone-LDDTestingOverview.pptx"

$picShape30 = $s30.Shapes.Item("Picture 1")
$picShape30.AlternativeText = "fig:  images/qr/overview_presentation.png"
